# Core64 LED Matrix V0.3 release to manufacturing
# - Update the "Package" (column G) footprint for all D1-D64 LED rows
#   from the old WS2813 PLCC6 footprint to the new "LED-6_5050" footprint.
# - Update the active cell selection left behind by the editor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 10 through 73 hold designators D1..D64 - update their Package (column G)
$pkgRange = $ws.Range("G10:G73")
$pkgRange.Value = "LED-6_5050"

# Match the formatting change applied alongside the new package name
# (an explicit Calibri font record, keeping the existing centered alignment)
$pkgRange.Font.Name = "Calibri"
$pkgRange.HorizontalAlignment = -4108  # xlCenter

# Leave the selection where the editor left it before saving
[void]$ws.Range("G8").Select()
